$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 becomes a real number instead of an inline/text string
$ws.Range("C2").Value = 6206872005

# New row 3 with contact data (all plain text strings)
$ws.Range("A3").Value = "jkljkl;lk"
$ws.Range("B3").Value = "isha@gmail.com"
$ws.Range("C3").Value = "jjor4455654"
$ws.Range("D3").Value = "sdkdklflkrk"
$ws.Range("E3").Value = "jlkjlkk;k;"
